$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row with "ПР09" (1-based row 20): first two date cells are empty and
# need "03.11" and "07.11" respectively.
$row1 = $t.Rows.Item(20)
$row1.Cells.Item(2).Range.Text = "03.11"
$row1.Cells.Item(3).Range.Text = "07.11"

# Row with "Л08" (1-based row 22): first date cell is empty and needs "04.11".
$row2 = $t.Rows.Item(22)
$row2.Cells.Item(2).Range.Text = "04.11"
